$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("D7").Value = "Yes"
[void]$ws.Range("E6").Select()
